# Generate Report for Handback
#
# For the "bf69e483-e797-4d6d-b849-7bd586fe366a" row (row 7) on both the
# zh-cn and de-de sheets, the handback run discovered that the handback
# file that came back is not the latest version. The report now fills in
# the previously-empty "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for that row, and
# widens the Error Detail column (P) so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb9276589b4eb12830f7e33d21f3b227203fa17e/e2e/bf69e483-e797-4d6d-b849-7bd586fe366a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab6d4a8287c937f95236ea70f9f4250139c94596/e2e/bf69e483-e797-4d6d-b849-7bd586fe366a.md."

# ColumnWidth is expressed in "characters"; the saved OOXML <col width>
# always carries ~0.8333333333333334 more than the ColumnWidth value that
# was set (Excel's column-width padding), so back that constant off to
# land on an OOXML width of exactly 40.
$targetColWidth = 40 - 0.8333333333333334

function Update-HandbackRow {
    param(
        [string]$SheetName,
        [string]$TargetFileUrl,
        [string]$HandbackFile,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File
    $ws.Range("I7").Value = "bf69e483-e797-4d6d-b849-7bd586fe366a.md"
    $ws.Hyperlinks.Add($ws.Range("I7"), $TargetFileUrl, "", "", "bf69e483-e797-4d6d-b849-7bd586fe366a.md")

    # Latest Handback File
    $ws.Range("J7").Value = $HandbackFile

    # Latest Handback DateTime
    $ws.Range("K7").Value = $HandbackDateTime

    # Error Detail
    $ws.Range("P7").Value = $errorDetail

    # Widen the Error Detail column so the message is readable.
    $ws.Columns.Item(16).ColumnWidth = $targetColWidth
}

Update-HandbackRow `
    "zh-cn" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bb9276589b4eb12830f7e33d21f3b227203fa17e/e2e/bf69e483-e797-4d6d-b849-7bd586fe366a.md" `
    "bf69e483-e797-4d6d-b849-7bd586fe366a.635b5f1cc06da6e57205996e198fc6da33821d9d.zh-cn.xlf" `
    "2016-08-28 12:43:12"

Update-HandbackRow `
    "de-de" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bb9276589b4eb12830f7e33d21f3b227203fa17e/e2e/bf69e483-e797-4d6d-b849-7bd586fe366a.md" `
    "bf69e483-e797-4d6d-b849-7bd586fe366a.635b5f1cc06da6e57205996e198fc6da33821d9d.de-de.xlf" `
    "2016-08-28 12:43:18"
